$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 111.655174
$ws.Cells.Item(33, 9).Value = 121
$ws.Cells.Item(33, 10).Value = 53.25
$ws.Cells.Item(33, 11).Value = 121
$ws.Cells.Item(33, 12).Value = 53.25
$ws.Cells.Item(33, 13).Value = 108
$ws.Cells.Item(33, 14).Value = -511.25
$ws.Cells.Item(100, 8).Value = 2946.3333
$ws.Cells.Item(100, 9).Value = 2952.125
$ws.Cells.Item(100, 10).Value = 2900
$ws.Cells.Item(100, 11).Value = 2952.125
$ws.Cells.Item(100, 12).Value = 2900
$ws.Cells.Item(100, 13).Value = -2411.125
$ws.Cells.Item(100, 14).Value = -3982
$ws.Cells.Item(112, 8).Value = 1367.1143
$ws.Cells.Item(112, 9).Value = 1379.8
$ws.Cells.Item(112, 10).Value = 1365
$ws.Cells.Item(112, 11).Value = 4139.4
$ws.Cells.Item(112, 12).Value = 4095
$ws.Cells.Item(112, 13).Value = -3031.4
$ws.Cells.Item(112, 14).Value = -6311
$ws.Cells.Item(116, 8).Value = 4449
$ws.Cells.Item(116, 9).Value = 5283.1665
$ws.Cells.Item(116, 10).Value = 3448
$ws.Cells.Item(116, 11).Value = 5283.1665
$ws.Cells.Item(116, 12).Value = 3448
$ws.Cells.Item(116, 13).Value = -1841.1665
$ws.Cells.Item(116, 14).Value = -10332
$ws.Cells.Item(123, 8).Value = 28000
$ws.Cells.Item(123, 10).Value = 28000
$ws.Cells.Item(123, 12).Value = 28000
$ws.Cells.Item(123, 14).Value = -37800
$ws.Cells.Item(126, 8).Value = 38880
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 38880
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 38880
$ws.Cells.Item(126, 13).ClearContents()
$ws.Cells.Item(126, 14).Value = -48760

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4583811.5
$ws.Cells.Item(32, 9).Value = 5147325.5
$ws.Cells.Item(32, 10).Value = 5262.5
$ws.Cells.Item(32, 11).Value = 5147325.5
$ws.Cells.Item(32, 12).Value = 5262.5
$ws.Cells.Item(32, 13).Value = -5147038.5
$ws.Cells.Item(32, 14).Value = -5836.5
$ws.Cells.Item(61, 8).Value = 1365.2
$ws.Cells.Item(61, 9).Value = 1000
$ws.Cells.Item(61, 10).Value = 1405.7778
$ws.Cells.Item(61, 11).Value = 1000
$ws.Cells.Item(61, 12).Value = 1405.7778
$ws.Cells.Item(61, 13).Value = -788
$ws.Cells.Item(61, 14).Value = -1829.7778
$ws.Cells.Item(74, 8).Value = 933.1667
$ws.Cells.Item(74, 9).Value = 979.8
$ws.Cells.Item(74, 10).Value = 700
$ws.Cells.Item(74, 11).Value = 979.8
$ws.Cells.Item(74, 12).Value = 700
$ws.Cells.Item(74, 13).Value = -105.8
$ws.Cells.Item(74, 14).Value = -2448
$ws.Cells.Item(77, 8).Value = 933.1667
$ws.Cells.Item(77, 9).Value = 979.8
$ws.Cells.Item(77, 10).Value = 700
$ws.Cells.Item(77, 11).Value = 4899
$ws.Cells.Item(77, 12).Value = 3500
$ws.Cells.Item(77, 13).Value = -531
$ws.Cells.Item(77, 14).Value = -12236
$ws.Cells.Item(88, 8).Value = 3470.6
$ws.Cells.Item(88, 9).Value = 3088.25
$ws.Cells.Item(88, 10).Value = 5000
$ws.Cells.Item(88, 11).Value = 3088.25
$ws.Cells.Item(88, 12).Value = 5000
$ws.Cells.Item(88, 13).Value = -2682.25
$ws.Cells.Item(88, 14).Value = -5812
$ws.Cells.Item(91, 8).Value = 3470.6
$ws.Cells.Item(91, 9).Value = 3088.25
$ws.Cells.Item(91, 10).Value = 5000
$ws.Cells.Item(91, 11).Value = 3088.25
$ws.Cells.Item(91, 12).Value = 5000
$ws.Cells.Item(91, 13).Value = -1684.25
$ws.Cells.Item(91, 14).Value = -7808
$ws.Cells.Item(136, 8).Value = 1365.2
$ws.Cells.Item(136, 9).Value = 1000
$ws.Cells.Item(136, 10).Value = 1405.7778
$ws.Cells.Item(136, 11).Value = 3000
$ws.Cells.Item(136, 12).Value = 4217.3334
$ws.Cells.Item(136, 13).Value = -450
$ws.Cells.Item(136, 14).Value = -9317.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 94014.95
$ws.Cells.Item(134, 9).Value = 3469.8125
$ws.Cells.Item(134, 10).Value = 335468.66
$ws.Cells.Item(134, 11).Value = 10409.4375
$ws.Cells.Item(134, 12).Value = 1006405.98
$ws.Cells.Item(134, 13).Value = -7874.4375
$ws.Cells.Item(134, 14).Value = -1011475.98

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 2581.081
$ws.Cells.Item(99, 9).Value = 2212.5
$ws.Cells.Item(99, 10).Value = 3261.5386
$ws.Cells.Item(99, 11).Value = 2212.5
$ws.Cells.Item(99, 12).Value = 3261.5386
$ws.Cells.Item(99, 13).Value = -714.5
$ws.Cells.Item(99, 14).Value = -6257.5386
$ws.Cells.Item(126, 8).Value = 2581.081
$ws.Cells.Item(126, 9).Value = 2212.5
$ws.Cells.Item(126, 10).Value = 3261.5386
$ws.Cells.Item(126, 11).Value = 6637.5
$ws.Cells.Item(126, 12).Value = 9784.6158
$ws.Cells.Item(126, 13).Value = -4167.5
$ws.Cells.Item(126, 14).Value = -14724.6158

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 12).Value = 0
$ws.Cells.Item(54, 14).ClearContents()
$ws.Cells.Item(59, 8).Value = 3166.3333
$ws.Cells.Item(59, 9).Value = 500
$ws.Cells.Item(59, 10).Value = 4499.5
$ws.Cells.Item(59, 11).Value = 1500
$ws.Cells.Item(59, 12).Value = 13498.5
$ws.Cells.Item(59, 13).Value = -960
$ws.Cells.Item(59, 14).Value = -14578.5
$ws.Cells.Item(64, 8).Value = 126213.875
$ws.Cells.Item(64, 9).Value = 428
$ws.Cells.Item(64, 10).Value = 251999.75
$ws.Cells.Item(64, 11).Value = 1284
$ws.Cells.Item(64, 12).Value = 755999.25
$ws.Cells.Item(64, 13).Value = -1014
$ws.Cells.Item(64, 14).Value = -756539.25
$ws.Cells.Item(67, 8).Value = 126213.875
$ws.Cells.Item(67, 9).Value = 428
$ws.Cells.Item(67, 10).Value = 251999.75
$ws.Cells.Item(67, 11).Value = 1284
$ws.Cells.Item(67, 12).Value = 755999.25
$ws.Cells.Item(67, 13).Value = -348
$ws.Cells.Item(67, 14).Value = -757871.25
$ws.Cells.Item(131, 8).Value = 7369100
$ws.Cells.Item(131, 9).Value = 50100396
$ws.Cells.Item(131, 10).Value = 1634.7069
$ws.Cells.Item(131, 11).Value = 150301188
$ws.Cells.Item(131, 12).Value = 4904.120699999999
$ws.Cells.Item(131, 13).Value = -150296148
$ws.Cells.Item(131, 14).Value = -14984.1207

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(123, 8).Value = 17299.834
$ws.Cells.Item(123, 10).Value = 17299.834
$ws.Cells.Item(123, 12).Value = 17299.834
$ws.Cells.Item(123, 14).Value = -22199.834

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2563.375
$ws.Cells.Item(7, 9).Value = 2501
$ws.Cells.Item(7, 10).Value = 3000
$ws.Cells.Item(7, 11).Value = 2501
$ws.Cells.Item(7, 12).Value = 3000
$ws.Cells.Item(7, 13).Value = -2389
$ws.Cells.Item(7, 14).Value = -3224
$ws.Cells.Item(40, 8).Value = 1264073.8
$ws.Cells.Item(40, 9).Value = 2021492
$ws.Cells.Item(40, 10).Value = 1710
$ws.Cells.Item(40, 11).Value = 2021492
$ws.Cells.Item(40, 12).Value = 1710
$ws.Cells.Item(40, 13).Value = -2021356
$ws.Cells.Item(40, 14).Value = -1982
$ws.Cells.Item(43, 8).Value = 253750
$ws.Cells.Item(43, 10).Value = 5000
$ws.Cells.Item(43, 12).Value = 5000
$ws.Cells.Item(43, 14).Value = -5386
$ws.Cells.Item(122, 8).Value = 10038.308
$ws.Cells.Item(122, 9).Value = 17083.166
$ws.Cells.Item(122, 10).Value = 3999.8572
$ws.Cells.Item(122, 11).Value = 51249.49800000001
$ws.Cells.Item(122, 12).Value = 11999.5716
$ws.Cells.Item(122, 13).Value = -48799.49800000001
$ws.Cells.Item(122, 14).Value = -16899.5716
$ws.Cells.Item(126, 8).Value = 2563.375
$ws.Cells.Item(126, 9).Value = 2501
$ws.Cells.Item(126, 10).Value = 3000
$ws.Cells.Item(126, 11).Value = 7503
$ws.Cells.Item(126, 12).Value = 9000
$ws.Cells.Item(126, 13).Value = -5033
$ws.Cells.Item(126, 14).Value = -13940
$ws.Cells.Item(136, 8).Value = 1967.1818
$ws.Cells.Item(136, 9).Value = 1334.6342
$ws.Cells.Item(136, 10).Value = 3819.6428
$ws.Cells.Item(136, 11).Value = 4003.9026
$ws.Cells.Item(136, 12).Value = 11458.9284
$ws.Cells.Item(136, 13).Value = -1453.9026
$ws.Cells.Item(136, 14).Value = -16558.9284

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1651.6666
$ws.Cells.Item(122, 9).Value = 1000
$ws.Cells.Item(122, 10).Value = 2303.3333
$ws.Cells.Item(122, 11).Value = 3000
$ws.Cells.Item(122, 12).Value = 6909.999899999999
$ws.Cells.Item(122, 13).Value = -550
$ws.Cells.Item(122, 14).Value = -11809.9999
$ws.Cells.Item(126, 8).Value = 685.7143
$ws.Cells.Item(126, 9).Value = 560
$ws.Cells.Item(126, 10).Value = 1000
$ws.Cells.Item(126, 11).Value = 1680
$ws.Cells.Item(126, 12).Value = 3000
$ws.Cells.Item(126, 13).Value = 790
$ws.Cells.Item(126, 14).Value = -7940
$ws.Cells.Item(136, 8).Value = 1668.2157
$ws.Cells.Item(136, 9).Value = 1714.579
$ws.Cells.Item(136, 11).Value = 5143.737
$ws.Cells.Item(136, 13).Value = -2593.737
